$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "254.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.48%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.249"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.79%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05850"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.48%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.712"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8665"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.91%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.042"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "21.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1409"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.41%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07176"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.18%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03183"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09232"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.56%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001544"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.25%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-94.08%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005804"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.78%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.65%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.03%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.47%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.45%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03449"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.85%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.34%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.538"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.99%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04152"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.91%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001228"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.17%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "15.88%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001201"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.18%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "1.28%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03803"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.41%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005785"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.77%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1099"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.66%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002350"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.89%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009691"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.48%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005239"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.32%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.17%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.09301"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "31.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002153"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-12.69%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.17%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
